$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 11.83567917099364
$ws.Range("C2").Value = 4.347754598935256
$ws.Range("D2").Value = 8.94793479154678
$ws.Range("E2").Value = 16.26076764462828
$ws.Range("F2").Value = 48.67805088739635
$ws.Range("H2").Value = 7.344005520526261
$ws.Range("K2").Value = 11.26397376326657
$ws.Range("B3").Value = 11.67813693589333
$ws.Range("C3").Value = 4.214988554051165
$ws.Range("D3").Value = 8.780091474645404
$ws.Range("E3").Value = 15.35067858270836
$ws.Range("F3").Value = 47.30757573533351
$ws.Range("H3").Value = 7.344005520526261
$ws.Range("K3").Value = 11.19011341761214
$ws.Range("B4").Value = 11.58636506302979
$ws.Range("C4").Value = 4.135341705279252
$ws.Range("D4").Value = 8.674561699303432
$ws.Range("E4").Value = 14.76924474061026
$ws.Range("F4").Value = 46.44854087176764
$ws.Range("H4").Value = 7.344005520526261
$ws.Range("K4").Value = 11.15011589207219
$ws.Range("B5").Value = 11.55026953716691
$ws.Range("C5").Value = 4.103436937274478
$ws.Range("D5").Value = 8.630961405431387
$ws.Range("E5").Value = 14.5268904505529
$ws.Range("F5").Value = 46.0944697873051
$ws.Range("H5").Value = 7.344005520526261
$ws.Range("K5").Value = 11.13517871219906
$ws.Range("B6").Value = 11.54435609469581
$ws.Range("C6").Value = 4.098174906308482
$ws.Range("D6").Value = 8.623686385272302
$ws.Range("E6").Value = 14.48632975882478
$ws.Range("F6").Value = 46.03544694542007
$ws.Range("H6").Value = 7.344005520526261
$ws.Range("K6").Value = 11.1327810834921
$ws.Range("B7").Value = 11.58587292335374
$ws.Range("C7").Value = 4.134909083424325
$ws.Range("D7").Value = 8.673976068866661
$ws.Range("E7").Value = 14.76599778755948
$ws.Range("F7").Value = 46.44378141327917
$ws.Range("H7").Value = 7.344005520526261
$ws.Range("K7").Value = 11.14990891000375
$ws.Range("B8").Value = 11.78036335212127
$ws.Range("C8").Value = 4.301636603089055
$ws.Range("D8").Value = 8.890590948904643
$ws.Range("E8").Value = 15.95181725065423
$ws.Range("F8").Value = 48.20939132980867
$ws.Range("H8").Value = 7.344005520526261
$ws.Range("K8").Value = 11.23740628242696
$ws.Range("B9").Value = 12.19852778953464
$ws.Range("C9").Value = 4.639975696130927
$ws.Range("D9").Value = 9.294669016677037
$ws.Range("E9").Value = 18.09520718372757
$ws.Range("F9").Value = 51.51592903586712
$ws.Range("H9").Value = 7.344005520526261
$ws.Range("K9").Value = 11.45060662414398
$ws.Range("B10").Value = 12.5245377967305
$ws.Range("C10").Value = 4.956420271253482
$ws.Range("D10").Value = 9.577567885519773
$ws.Range("E10").Value = 19.7183699563193
$ws.Range("F10").Value = 53.83016464184863
$ws.Range("H10").Value = 7.344005520526261
$ws.Range("K10").Value = 11.63127218317543
$ws.Range("B11").Value = 12.67612780060946
$ws.Range("C11").Value = 5.13513772057759
$ws.Range("D11").Value = 9.702968576233319
$ws.Range("E11").Value = 20.41597673032961
$ws.Range("F11").Value = 54.85430601257142
$ws.Range("H11").Value = 7.344005520526261
$ws.Range("K11").Value = 11.71834083024936
$ws.Range("B12").Value = 12.73393651624357
$ws.Range("C12").Value = 5.201349327230122
$ws.Range("D12").Value = 9.749960865131113
$ws.Range("E12").Value = 20.67431822616799
$ws.Range("F12").Value = 55.23773712079586
$ws.Range("H12").Value = 7.344005520526261
$ws.Range("K12").Value = 11.75198229709537
$ws.Range("B13").Value = 12.72146943886117
$ws.Range("C13").Value = 5.187154467468088
$ws.Range("D13").Value = 9.739862534990964
$ws.Range("E13").Value = 20.61893833800332
$ws.Range("F13").Value = 55.15535765127839
$ws.Range("H13").Value = 7.344005520526261
$ws.Range("K13").Value = 11.744707708544
$ws.Range("B14").Value = 12.68087602963644
$ws.Range("C14").Value = 5.140614291093901
$ws.Range("D14").Value = 9.706844680352811
$ws.Range("E14").Value = 20.4373472131121
$ws.Range("F14").Value = 54.88594055689418
$ws.Range("H14").Value = 7.344005520526261
$ws.Range("K14").Value = 11.72109528440813
$ws.Range("B15").Value = 12.6560621221567
$ws.Range("C15").Value = 5.111916613627521
$ws.Range("D15").Value = 9.686555329300091
$ws.Range("E15").Value = 20.32535964171173
$ws.Range("F15").Value = 54.7203355560064
$ws.Range("H15").Value = 7.344005520526261
$ws.Range("K15").Value = 11.70671833722692
$ws.Range("B16").Value = 12.51469136335718
$ws.Range("C16").Value = 4.944535161009279
$ws.Range("D16").Value = 9.569304604059447
$ws.Range("E16").Value = 19.67196230152317
$ws.Range("F16").Value = 53.76263472220572
$ws.Range("H16").Value = 7.344005520526261
$ws.Range("K16").Value = 11.62567762878973
$ws.Range("B17").Value = 12.42875837701747
$ws.Range("C17").Value = 4.839232722154143
$ws.Range("D17").Value = 9.496516754917586
$ws.Range("E17").Value = 19.26070616827268
$ws.Range("F17").Value = 53.16758341137174
$ws.Range("H17").Value = 7.344005520526261
$ws.Range("K17").Value = 11.57719041949066
$ws.Range("B18").Value = 12.37964536786099
$ws.Range("C18").Value = 4.780592728263606
$ws.Range("D18").Value = 9.454342471934851
$ws.Range("E18").Value = 19.0203220969524
$ws.Range("F18").Value = 52.82264678763443
$ws.Range("H18").Value = 7.344005520526261
$ws.Range("K18").Value = 11.54976348566148
$ws.Range("B19").Value = 12.36307244055631
$ws.Range("C19").Value = 4.767846018382345
$ws.Range("D19").Value = 9.440010635263357
$ws.Range("E19").Value = 18.9382714620453
$ws.Range("F19").Value = 52.70540587398288
$ws.Range("H19").Value = 7.344005520526261
$ws.Range("K19").Value = 11.54055741972867
$ws.Range("B20").Value = 12.43787414973991
$ws.Range("C20").Value = 4.850542222962359
$ws.Range("D20").Value = 9.504297244345784
$ws.Range("E20").Value = 19.30488228760303
$ws.Range("F20").Value = 53.23120696085179
$ws.Range("H20").Value = 7.344005520526261
$ws.Range("K20").Value = 11.58230441343479
$ws.Range("B21").Value = 12.69278884607351
$ws.Range("C21").Value = 5.154323946790879
$ws.Range("D21").Value = 9.716556398066171
$ws.Range("E21").Value = 20.49084278024596
$ws.Range("F21").Value = 54.96519605525668
$ws.Range("H21").Value = 7.344005520526261
$ws.Range("K21").Value = 11.72801289047393
$ws.Range("B22").Value = 12.86171943648857
$ws.Range("C22").Value = 5.344329282107955
$ws.Range("D22").Value = 9.852390507646664
$ws.Range("E22").Value = 21.23198899797753
$ws.Range("F22").Value = 56.07276731222988
$ws.Range("H22").Value = 7.344005520526261
$ws.Range("K22").Value = 11.82713396140721
$ws.Range("B23").Value = 12.77136705321089
$ws.Range("C23").Value = 5.243697497079832
$ws.Range("D23").Value = 9.780164245486588
$ws.Range("E23").Value = 20.83951939242089
$ws.Range("F23").Value = 55.48407004626327
$ws.Range("H23").Value = 7.344005520526261
$ws.Range("K23").Value = 11.77388593053781
$ws.Range("B24").Value = 12.43375199731877
$ws.Range("C24").Value = 4.845432291676664
$ws.Range("D24").Value = 9.500780703808122
$ws.Range("E24").Value = 19.28492256224138
$ws.Range("F24").Value = 53.20245155938685
$ws.Range("H24").Value = 7.344005520526261
$ws.Range("K24").Value = 11.57999097582579
$ws.Range("B25").Value = 12.08183514572718
$ws.Range("C25").Value = 4.547697622136201
$ws.Range("D25").Value = 9.187716290592022
$ws.Range("E25").Value = 17.53147427128374
$ws.Range("F25").Value = 50.64028923727951
$ws.Range("H25").Value = 7.344005520526261
$ws.Range("K25").Value = 11.38860331065441
